$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'66.954.57"
$ws.Range("E2").Value = "  +2.36%  "

$ws.Range("D3").Formula = "'3.105.48"
$ws.Range("E3").Value = "  +5.25%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Formula = "'579.98"
$ws.Range("E5").Value = "  +1.81%  "

$ws.Range("D6").Formula = "'173.00"
$ws.Range("E6").Value = "  +8.02%  "

$ws.Range("D7").Formula = "'1.00"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Formula = "'3.101.24"
$ws.Range("E8").Value = "  +5.28%  "

$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("E11").Value = "  +3.99%  "

$ws.Range("D12").Formula = "'0.482"
$ws.Range("E12").Value = "  +4.91%  "

$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").Formula = "'37.39"
$ws.Range("E14").Value = "  +8.33%  "

$ws.Range("D16").Formula = "'3.617.88"
$ws.Range("E16").Value = "  +5.18%  "

$ws.Range("D17").Formula = "'66.910.67"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Formula = "'7.20"
$ws.Range("E18").Value = "  +2.62%  "

$ws.Range("D19").Formula = "'3.103.60"
$ws.Range("E19").Value = "  +5.18%  "

$ws.Range("D20").Formula = "'16.20"
$ws.Range("E20").Value = "  +3.41%  "

$ws.Range("D21").Formula = "'483.08"
$ws.Range("E21").Value = "  +8.49%  "

$ws.Range("D22").Formula = "'0.715"
$ws.Range("E22").Value = "  +2.89%  "

$ws.Range("D23").Formula = "'7.53"
$ws.Range("E23").Value = "  +3.20%  "

$ws.Range("D24").Formula = "'84.10"
$ws.Range("E24").Value = "  +2.23%  "

$ws.Range("D25").Formula = "'2.35"
$ws.Range("E25").Value = "  +5.41%  "

$ws.Range("D26").Formula = "'13.01"
$ws.Range("E26").Value = "  +6.65%  "

$ws.Range("D27").Formula = "'10.04"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Formula = "'7.98"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("D31").Formula = "'2.69"
$ws.Range("E31").Value = "  +3.93%  "

$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("E33").Value = "  +5.59%  "

$ws.Range("E34").Value = "  +2.26%  "

$ws.Range("D35").Formula = "'0.999"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  +3.43%  "

$ws.Range("D37").Formula = "'5.89"
$ws.Range("E37").Value = "  +2.82%  "

$ws.Range("D38").Formula = "'48.24"
$ws.Range("E38").Value = "  +7.52%  "

$ws.Range("D39").Formula = "'2.13"
$ws.Range("E39").Value = "  +8.37%  "

$ws.Range("D40").Formula = "'50.21"
$ws.Range("E40").Value = "  +2.27%  "

$ws.Range("E41").Value = "  +5.09%  "

$ws.Range("D43").Formula = "'8.66"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("E44").Value = "  -1.64%  "

$ws.Range("D45").Formula = "'2.825.89"
$ws.Range("E45").Value = "  +5.47%  "

$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("D47").Formula = "'380.27"
$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("D48").Formula = "'135.33"
$ws.Range("E48").Value = "  +1.36%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").Formula = "'24.87"
$ws.Range("E50").Value = "  +5.17%  "

$ws.Range("E51").Value = "  +2.67%  "
